# TouchShield-Pin-Assignment.xlsx edit script
# Renames the "per channel" detect-pin naming to "per sensor" naming:
#   DETECT_CH##                              -> DETECT-SENSOR##
#   SHIELD-DETECT                            -> SHIELD-DETECT-OUTPUT
#   Digital output for sense detect on channel ##  -> Digital output for detection on sensor ##
# Also adds a note to the USB DM/DP rows, and restores the "Pin Assignment"
# sheet/cell selection that was active when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin Assignment")

# Rows (in sheet order) that hold the DETECT_CH00..DETECT_CH31 pins, in
# numeric order - the Nth row in this list is channel/sensor number N.
$detectRows = @(25,26,27,28,29,30,31,32,33,34,35,38,41,42,43,44,45,46,47,48,49,50,122,123,124,125,126,127,128,129,130,131)

$sensorNum = 0
foreach ($r in $detectRows) {
    $numStr = "{0:D2}" -f $sensorNum
    $ws.Cells.Item($r, 4).Value = "DETECT-SENSOR$numStr"
    $ws.Cells.Item($r, 5).Value = "SHIELD-DETECT-OUTPUT"
    $ws.Cells.Item($r, 6).Value = "Digital output for detection on sensor $numStr"
    $sensorNum = $sensorNum + 1
}

# DM / DP (USB data) pins: note the possible use for USB firmware upgrade
$ws.Cells.Item(52, 10).Value = "??? Perhaps use for USB firmware upgrade ???"
$ws.Cells.Item(53, 10).Value = "??? Perhaps use for USB firmware upgrade ???"

# Column D/E got a bit wider once "DETECT-SENSOR00".."DETECT-SENSOR31" and
# "SHIELD-DETECT-OUTPUT" (longer than the old names) are in use.
$ws.Columns.Item(4).ColumnWidth = 24
$ws.Columns.Item(5).ColumnWidth = 26

# Restore view: "Pin Assignment" tab active again, selection moved to A57.
$ws.Activate()
$ws.Range("A57").Select()
